$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.1050873333333333"
$ws.Range("H2").Value = [double]"0.315262"
$ws.Range("I2").Value = [double]"0.03031434174852429"
$ws.Range("J2").Value = [double]"0.03031434174852429"
$ws.Range("M2").Value = [double]"28.689524"
$ws.Range("N2").Value = [double]"86.06857199999999"
$ws.Range("O2").Value = [double]"0.2394085694101769"
$ws.Range("P2").Value = [double]"0.2394085694101769"
$ws.Range("Q2").Value = [double]"3.014905571762666"
$ws.Range("R2").Value = [double]"27.13415014586399"
$ws.Range("S2").Value = [double]"0.007257513190625403"
$ws.Range("T2").Value = [double]"0.007257513190625403"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.1050873333333333"
$ws.Range("H3").Value = [double]"0.315262"
$ws.Range("I3").Value = [double]"0.03031434174852429"
$ws.Range("J3").Value = [double]"0.03031434174852429"
$ws.Range("O3").Value = [double]"0.5212694246546397"
$ws.Range("P3").Value = [double]"0.5212694246546395"
$ws.Range("Q3").Value = [double]"6.564418711713778"
$ws.Range("R3").Value = [double]"59.079768405424"
$ws.Range("S3").Value = [double]"0.01580193948203738"
$ws.Range("T3").Value = [double]"0.01580193948203738"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.1050873333333333"
$ws.Range("H4").Value = [double]"0.315262"
$ws.Range("I4").Value = [double]"0.03031434174852429"
$ws.Range("J4").Value = [double]"0.03031434174852429"
$ws.Range("M4").Value = [double]"28.525746"
$ws.Range("N4").Value = [double]"85.57723799999999"
$ws.Range("O4").Value = [double]"0.2380418734454457"
$ws.Range("P4").Value = [double]"0.2380418734454457"
$ws.Range("Q4").Value = [double]"2.997694578483999"
$ws.Range("R4").Value = [double]"26.979251206356"
$ws.Range("S4").Value = [double]"0.007216082702084211"
$ws.Range("T4").Value = [double]"0.007216082702084211"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.1050873333333333"
$ws.Range("H5").Value = [double]"0.315262"
$ws.Range("I5").Value = [double]"0.03031434174852429"
$ws.Range("J5").Value = [double]"0.03031434174852429"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.1534046666666667"
$ws.Range("N5").Value = [double]"0.460214"
$ws.Range("O5").Value = [double]"0.001280132489737778"
$ws.Range("P5").Value = [double]"0.001280132489737778"
$ws.Range("Q5").Value = [double]"0.01612088734088889"
$ws.Range("R5").Value = [double]"0.145087986068"
$ws.Range("S5").Value = [double]"3.880637377730026E-05"
$ws.Range("T5").Value = [double]"3.880637377730026E-05"
$ws.Range("I6").Value = [double]"0.6354599969768544"
$ws.Range("J6").Value = [double]"0.6354599969768545"
$ws.Range("M6").Value = [double]"28.689524"
$ws.Range("N6").Value = [double]"86.06857199999999"
$ws.Range("O6").Value = [double]"0.2394085694101769"
$ws.Range("P6").Value = [double]"0.2394085694101769"
$ws.Range("Q6").Value = [double]"63.19952125007198"
$ws.Range("R6").Value = [double]"568.7956912506479"
$ws.Range("S6").Value = [double]"0.1521345687936241"
$ws.Range("T6").Value = [double]"0.1521345687936241"
$ws.Range("I7").Value = [double]"0.6354599969768544"
$ws.Range("J7").Value = [double]"0.6354599969768545"
$ws.Range("O7").Value = [double]"0.5212694246546397"
$ws.Range("P7").Value = [double]"0.5212694246546395"
$ws.Range("S7").Value = [double]"0.331245867015164"
$ws.Range("T7").Value = [double]"0.331245867015164"
$ws.Range("I8").Value = [double]"0.6354599969768544"
$ws.Range("J8").Value = [double]"0.6354599969768545"
$ws.Range("M8").Value = [double]"28.525746"
$ws.Range("N8").Value = [double]"85.57723799999999"
$ws.Range("O8").Value = [double]"0.2380418734454457"
$ws.Range("P8").Value = [double]"0.2380418734454457"
$ws.Range("Q8").Value = [double]"62.83873829698798"
$ws.Range("R8").Value = [double]"565.5486446728919"
$ws.Range("S8").Value = [double]"0.1512660881800077"
$ws.Range("T8").Value = [double]"0.1512660881800077"
$ws.Range("I9").Value = [double]"0.6354599969768544"
$ws.Range("J9").Value = [double]"0.6354599969768545"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.1534046666666667"
$ws.Range("N9").Value = [double]"0.460214"
$ws.Range("O9").Value = [double]"0.001280132489737778"
$ws.Range("P9").Value = [double]"0.001280132489737778"
$ws.Range("Q9").Value = [double]"0.3379317652973333"
$ws.Range("R9").Value = [double]"3.041385887676"
$ws.Range("S9").Value = [double]"0.0008134729880587414"
$ws.Range("T9").Value = [double]"0.0008134729880587415"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"0.9666886666666668"
$ws.Range("H10").Value = [double]"2.900066"
$ws.Range("I10").Value = [double]"0.2788588279503266"
$ws.Range("J10").Value = [double]"0.2788588279503266"
$ws.Range("M10").Value = [double]"28.689524"
$ws.Range("N10").Value = [double]"86.06857199999999"
$ws.Range("O10").Value = [double]"0.2394085694101769"
$ws.Range("P10").Value = [double]"0.2394085694101769"
$ws.Range("Q10").Value = [double]"27.73383770286133"
$ws.Range("R10").Value = [double]"249.604539325752"
$ws.Range("S10").Value = [double]"0.06676119306698636"
$ws.Range("T10").Value = [double]"0.06676119306698636"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"0.9666886666666668"
$ws.Range("H11").Value = [double]"2.900066"
$ws.Range("I11").Value = [double]"0.2788588279503266"
$ws.Range("J11").Value = [double]"0.2788588279503266"
$ws.Range("O11").Value = [double]"0.5212694246546397"
$ws.Range("P11").Value = [double]"0.5212694246546395"
$ws.Range("Q11").Value = [double]"60.3854810145369"
$ws.Range("R11").Value = [double]"543.4693291308321"
$ws.Range("S11").Value = [double]"0.1453605808055339"
$ws.Range("T11").Value = [double]"0.1453605808055339"
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"0.9666886666666668"
$ws.Range("H12").Value = [double]"2.900066"
$ws.Range("I12").Value = [double]"0.2788588279503266"
$ws.Range("J12").Value = [double]"0.2788588279503266"
$ws.Range("M12").Value = [double]"28.525746"
$ws.Range("N12").Value = [double]"85.57723799999999"
$ws.Range("O12").Value = [double]"0.2380418734454457"
$ws.Range("P12").Value = [double]"0.2380418734454457"
$ws.Range("Q12").Value = [double]"27.575515366412"
$ws.Range("R12").Value = [double]"248.179638297708"
$ws.Range("S12").Value = [double]"0.06638007783209697"
$ws.Range("T12").Value = [double]"0.06638007783209697"
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"0.9666886666666668"
$ws.Range("H13").Value = [double]"2.900066"
$ws.Range("I13").Value = [double]"0.2788588279503266"
$ws.Range("J13").Value = [double]"0.2788588279503266"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"0.1534046666666667"
$ws.Range("N13").Value = [double]"0.460214"
$ws.Range("O13").Value = [double]"0.001280132489737778"
$ws.Range("P13").Value = [double]"0.001280132489737778"
$ws.Range("Q13").Value = [double]"0.1482945526804444"
$ws.Range("R13").Value = [double]"1.334650974124"
$ws.Range("S13").Value = [double]"0.0003569762457094102"
$ws.Range("T13").Value = [double]"0.0003569762457094102"
$ws.Range("G14").Value = [double]"0.191934"
$ws.Range("H14").Value = [double]"0.575802"
$ws.Range("I14").Value = [double]"0.05536683332429468"
$ws.Range("J14").Value = [double]"0.05536683332429467"
$ws.Range("M14").Value = [double]"28.689524"
$ws.Range("N14").Value = [double]"86.06857199999999"
$ws.Range("O14").Value = [double]"0.2394085694101769"
$ws.Range("P14").Value = [double]"0.2394085694101769"
$ws.Range("Q14").Value = [double]"5.506495099416"
$ws.Range("R14").Value = [double]"49.558455894744"
$ws.Range("S14").Value = [double]"0.0132552943589411"
$ws.Range("T14").Value = [double]"0.0132552943589411"
$ws.Range("G15").Value = [double]"0.191934"
$ws.Range("H15").Value = [double]"0.575802"
$ws.Range("I15").Value = [double]"0.05536683332429468"
$ws.Range("J15").Value = [double]"0.05536683332429467"
$ws.Range("O15").Value = [double]"0.5212694246546397"
$ws.Range("P15").Value = [double]"0.5212694246546395"
$ws.Range("Q15").Value = [double]"11.989410151056"
$ws.Range("R15").Value = [double]"107.904691359504"
$ws.Range("S15").Value = [double]"0.02886103735190442"
$ws.Range("T15").Value = [double]"0.02886103735190441"
$ws.Range("G16").Value = [double]"0.191934"
$ws.Range("H16").Value = [double]"0.575802"
$ws.Range("I16").Value = [double]"0.05536683332429468"
$ws.Range("J16").Value = [double]"0.05536683332429467"
$ws.Range("M16").Value = [double]"28.525746"
$ws.Range("N16").Value = [double]"85.57723799999999"
$ws.Range("O16").Value = [double]"0.2380418734454457"
$ws.Range("P16").Value = [double]"0.2380418734454457"
$ws.Range("Q16").Value = [double]"5.475060532764"
$ws.Range("R16").Value = [double]"49.275544794876"
$ws.Range("S16").Value = [double]"0.01317962473125684"
$ws.Range("T16").Value = [double]"0.01317962473125684"
$ws.Range("G17").Value = [double]"0.191934"
$ws.Range("H17").Value = [double]"0.575802"
$ws.Range("I17").Value = [double]"0.05536683332429468"
$ws.Range("J17").Value = [double]"0.05536683332429467"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"0.1534046666666667"
$ws.Range("N17").Value = [double]"0.460214"
$ws.Range("O17").Value = [double]"0.001280132489737778"
$ws.Range("P17").Value = [double]"0.001280132489737778"
$ws.Range("Q17").Value = [double]"0.029443571292"
$ws.Range("R17").Value = [double]"0.264992141628"
$ws.Range("S17").Value = [double]"7.087688219232591E-05"
$ws.Range("T17").Value = [double]"7.087688219232589E-05"
